$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 120.28

# Row 3
$ws.Range("E3").Value = 294.84

# Row 4
$ws.Range("E4").Value = 80.16

# Row 5
$ws.Range("E5").Value = 1747.46

# Row 8
$ws.Range("D8").Value = "6012A10-1, 6012A10-1, 1929A18-1, 1918A10-1, 1911A10-1, 1973C00-1, 1972A00-1, 1917A10-1, 1929A18-1, 1932A10-1, 6001A10-1"
$ws.Range("E8").Value = 669.13
$ws.Range("F8").Value = 0.43090675

# Row 9
$ws.Range("E9").Value = 115.74

# Row 10
$ws.Range("E10").Value = 11.79

# Row 11
$ws.Range("E11").Value = 90.73499999999999

# Row 12
$ws.Range("E12").Value = 1305.28

# Row 13
$ws.Range("D13").Value = "6001A10-1, 6001A10-1"
$ws.Range("F13").Value = 0.1121

# Row 14
$ws.Range("E14").Value = 316.24

# Row 15
$ws.Range("E15").Value = 190.184

# Row 16
$ws.Range("E16").Value = 17.36

# Row 17
$ws.Range("E17").Value = 17.2

# Row 18
$ws.Range("E18").Value = 44.4

# Row 20
$ws.Range("E20").Value = 2.84

# Row 21
$ws.Range("D21").Value = "SB75A04-1, SB57A04-1, SB75A04-1, SB65A04-1"
$ws.Range("E21").Value = 13.184
$ws.Range("F21").Value = 0.01817235

# Row 22
$ws.Range("E22").Value = 12.954

# Row 23
$ws.Range("D23").Value = "SB82A04-1, SB57A04-1, SB57A04-1, SB82A04-1"
$ws.Range("E23").Value = 29.724
$ws.Range("F23").Value = 0.038664

# Row 24
$ws.Range("D24").Value = "SB51A04-1, SB82A04-1, SB57A04-1, SB82A04-1, SB51A04-1"
$ws.Range("E24").Value = 88.75
$ws.Range("F24").Value = 0.119088

# Row 25
$ws.Range("E25").Value = 12.974

# Row 26
$ws.Range("E26").Value = 24.1

# Row 27
$ws.Range("E27").Value = 109.884

# Row 28
$ws.Range("D28").Value = "SB62A04-1, SB51A04-1, SB49A04-1, SB50A04-1, SB27A04-1, SB27A04-1, SB37A04-1"
$ws.Range("E28").Value = 208.88
$ws.Range("F28").Value = 0.234938185

# Row 29
$ws.Range("E29").Value = 39.768

# Row 31
$ws.Range("E31").Value = 851.3699999999999

# Row 32
$ws.Range("E32").Value = 27.88

# Row 33
$ws.Range("E33").Value = 15.2

# Row 35
$ws.Range("E35").Value = 7.632

# Row 37
$ws.Range("D37").Value = "04DBA04-1, 06DHA04-1, 06DHA04-1, 06DJA04-1"
$ws.Range("E37").Value = 35.48
$ws.Range("F37").Value = 0.02962604

# Row 38
$ws.Range("E38").Value = 2010.925

# Row 39
$ws.Range("E39").Value = 21828.8

# Row 40
$ws.Range("E40").Value = 19877

# Row 41
$ws.Range("E41").Value = 12455.2

# Row 42
$ws.Range("E42").Value = 12.36

# Row 43
$ws.Range("D43").Value = "06CKA04-1, 4941A04-1, 0604A04-1, 0843A00-1, 4009A04-1, 06CKA04-1, 4009A04-1"
$ws.Range("E43").Value = 280.735
$ws.Range("F43").Value = 0.129596905

# Row 44
$ws.Range("E44").Value = 45.32

# Row 45
$ws.Range("D45").Value = "06C9A04-1, 1088A04-1, 0604A04-1, 04BNA04-1, 04BGA04-1, 06CDA04-1, 0619A04-1, 04BNA04-1, 01AHA01-1, 06ASA04-1, 06CEA04-1, 06CEA04-1, 04BNA04-1, 1088A04-1, 01AHA01-1, 06CEA04-1, 01AHA01-1, 04BNA04-1, 06C4A04-1, 04BNA04-1, 06C9A04-1"
$ws.Range("E45").Value = 5669.58
$ws.Range("F45").Value = 0.95389843

# Row 46
$ws.Range("D46").Value = "0619A04-1, 06CEA04-1, 06CEA04-1, 06ASA04-1, 06ASA04-1, 13A5A03-1, 01BWA03-1, 4017A04-1, 13FTA03-1, 13CBA03-1, 0601A04-1, 04C1A04-1"
$ws.Range("E46").Value = 1148.88
$ws.Range("F46").Value = 0.17302811

# Row 47
$ws.Range("E47").Value = 48.81

# Row 48
$ws.Range("E48").Value = 22.58

# Row 49
$ws.Range("E49").Value = 91.26000000000001

# Row 50
$ws.Range("D50").Value = "04BNA04-1, 04BNA04-1, 06ASA04-1, 04BNA04-1"
$ws.Range("E50").Value = 48.8
$ws.Range("F50").Value = 0.01843065

# Row 51
$ws.Range("D51").Value = "06CEA04-1, 06CEA04-1"
$ws.Range("F51").Value = 0.1404

# Row 52
$ws.Range("E52").Value = 766.0529999999999

# Row 53
$ws.Range("D53").Value = "0913B01-1, 01AHA01-1, 06ASA04-1, 06ASA04-1, 0843A00-1, 0924A00-1, 01AHA01-1"
$ws.Range("E53").Value = 379.47
$ws.Range("F53").Value = 0.0679448

# Row 54
$ws.Range("D54").Value = "1395A01-1, 06CDA04-1, 06ASA04-1, 06ASA04-1, 3143A08-1, 0843A00-1, 0649A04-1"
$ws.Range("E54").Value = 288.33
$ws.Range("F54").Value = 0.07165515

# Row 55
$ws.Range("D55").Value = "1395A01-1, 01AHA01-1, 06ASA04-1, 01AHA01-1"
$ws.Range("E55").Value = 95.68000000000001
$ws.Range("F55").Value = 0.0316975

# Row 57
$ws.Range("E57").Value = 76.12

# Row 58
$ws.Range("E58").Value = 127.2

# Row 59
$ws.Range("D59").Value = "04BNA04-1, 06CEA04-1, 06CEA04-1"
$ws.Range("E59").Value = 181.11
$ws.Range("F59").Value = 0.1059412

# Row 60
$ws.Range("D60").Value = "04BNA04-1, 06CEA04-1, 06CEA04-1, 06ASA04-1, 0913B01-1, 0956A00-1, 3143A08-1, 06CEA04-1, 06C6A04-1, 0928A00-1"
$ws.Range("E60").Value = 983.7
$ws.Range("F60").Value = 0.22209843

# Row 61
$ws.Range("E61").Value = 10.86

# Row 62
$ws.Range("D62").Value = "06CDA04-1, 06BLA04-1, 06ASA04-1, 06ASA04-1"
$ws.Range("E62").Value = 158.892
$ws.Range("F62").Value = 0.120428

# Row 63
$ws.Range("E63").Value = 53.40000000000001

# Row 64
$ws.Range("D64").Value = "0916A00-1, 0916A00-1, 3143A08-1, 04BNA04-1, 0601A04-1, 0619A04-1, 0619A04-1, 0916A00-1"
$ws.Range("E64").Value = 1436.224
$ws.Range("F64").Value = 0.392123634

# Row 65
$ws.Range("D65").Value = "0916A00-1, 0944A00-1, 04BNA04-1, 06C6A04-1, 0916A00-1, 0154A03-1, 01HJA03-1, 01H8A03-1, 01G4A03-1, 0125A03-1, 1436A01-1, 1306A01-1, 13FDA03-1, 2524A03-1, 2551A03-1"
$ws.Range("E65").Value = 2133.825
$ws.Range("F65").Value = 0.1721447

# Row 66
$ws.Range("D66").Value = "0649A04-1, 01AHA01-1, 04BNA04-1, 0928A00-1, 01AHA01-1, 04BNA04-1, 01AHA01-1"
$ws.Range("E66").Value = 2681.714
$ws.Range("F66").Value = 0.486864886

# Row 67
$ws.Range("E67").Value = 442.295

# Row 68
$ws.Range("E68").Value = 1127.39

# Row 69
$ws.Range("E69").Value = 1131.975

# Row 70
$ws.Range("E70").Value = 20.75

# Row 71
$ws.Range("E71").Value = 1947.918

# Row 72
$ws.Range("E72").Value = 50.64000000000001

# Row 84
$ws.Range("E84").Value = 0.4

# Row 85
$ws.Range("E85").Value = 30.04

# Row 86
$ws.Range("E86").Value = 32941.16800000001

# Row 87
$ws.Range("E87").Value = 786.27
